$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.110.06"
$ws.Range("E2").Value = "  -0.96%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.971.92"
$ws.Range("E3").Value = "  -1.26%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = "  +0.18%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'329.60"
$ws.Range("E5").Value = "  -0.51%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.24%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.4993"
$ws.Range("E7").Value = "  +1.10%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.4223"
$ws.Range("E8").Value = "  +1.01%  "

# Row 9 - OKB
$ws.Range("D9").Value = "'53.02"
$ws.Range("E9").Value = "  -0.82%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.09257"
$ws.Range("E10").Value = "  +5.46%  "

# Row 11 - Polygon
$ws.Range("D11").Value = "'1.102"
$ws.Range("E11").Value = "  -0.85%  "

# Row 12 - Solana
$ws.Range("D12").Value = "'22.90"
$ws.Range("E12").Value = "  -0.75%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.991.05"
$ws.Range("E13").Value = "  -9.66%  "

# Row 14 - Chainlink
$ws.Range("D14").Value = "'7.917"
$ws.Range("E14").Value = "  -1.88%  "

# Row 15 - Polkadot
$ws.Range("D15").Value = "'6.467"
$ws.Range("E15").Value = "  -0.01%  "

# Row 16 - BinanceUSD
$ws.Range("E16").Value = "  +0.16%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "'0.00001109"
$ws.Range("E17").Value = "  +0.49%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "'91.84"
$ws.Range("E18").Value = "  -4.41%  "

# Row 19 - TRON
$ws.Range("D19").Value = "'0.06737"
$ws.Range("E19").Value = "  +1.55%  "

# Row 20 - Avalanche
$ws.Range("E20").Value = "  -0.80%  "

# Row 21 - Dai
$ws.Range("D21").Value = "'1.008"
$ws.Range("E21").Value = "  +0.08%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'5.974"
$ws.Range("E22").Value = "  +0.26%  "

# Row 23 - WrappedBTC
$ws.Range("D23").Value = "29.123.62"
$ws.Range("E23").Value = "  -1.21%  "

# Row 24 - Cosmos
$ws.Range("D24").Value = "'11.96"
$ws.Range("E24").Value = "  +1.64%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "'2.265"
$ws.Range("E25").Value = "  -0.82%  "

# Row 26 - WrappedliquidstakedEther2.0
$ws.Range("D26").Value = "2.227.81"
$ws.Range("E26").Value = "  -4.68%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'20.72"
$ws.Range("E27").Value = "  +1.09%  "

# Row 28 - Monero
$ws.Range("D28").Value = "'155.65"
$ws.Range("E28").Value = "  -1.19%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "'6.318"
$ws.Range("E29").Value = "  -5.30%  "

# Row 30 - LidoDAOToken
$ws.Range("D30").Value = "'2.265"
$ws.Range("E30").Value = "  -3.31%  "

# Row 31 - BitcoinCash
$ws.Range("D31").Value = "'126.75"
$ws.Range("E31").Value = "  +0.05%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "'1.050"
$ws.Range("E32").Value = "  +0.43%  "

# Row 33 - Stellar
$ws.Range("D33").Value = "'0.09870"
$ws.Range("E33").Value = "  -0.37%  "

# Row 34 - ARBITRUM
$ws.Range("D34").Value = "'1.519"
$ws.Range("E34").Value = "  -2.03%  "

# Row 35 - Filecoin
$ws.Range("D35").Value = "'5.828"
$ws.Range("E35").Value = "  +0.42%  "

# Row 36 - HuobiToken
$ws.Range("D36").Value = "'3.719"
$ws.Range("E36").Value = "  -1.68%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "'0.02437"
$ws.Range("E37").Value = "  -0.25%  "

# Row 38 - TrustWalletToken
$ws.Range("D38").Value = "'1.318"
$ws.Range("E38").Value = "  +2.56%  "

# Row 39 - FraxShare -> Hedera (identity swap with row 40)
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.06407"
$ws.Range("E39").Value = "  +0.71%  "

# Row 40 - Hedera -> FraxShare (identity swap with row 39)
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'9.041"
$ws.Range("E40").Value = "  -5.95%  "

# Row 41 - TheSandbox
$ws.Range("D41").Value = "'0.6491"
$ws.Range("E41").Value = "  +0.18%  "

# Row 42 - Aptos
$ws.Range("D42").Value = "'11.49"
$ws.Range("E42").Value = "  -1.99%  "

# Row 43 - Algorand
$ws.Range("D43").Value = "'0.2004"
$ws.Range("E43").Value = "  -2.66%  "

# Row 44 - Frax
$ws.Range("E44").Value = "  +0.03%  "

# Row 45 - Decentraland -> WEMIXTOKEN (identity swap with row 46)
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.378"
$ws.Range("E45").Value = "  +9.80%  "

# Row 46 - WEMIXTOKEN -> Decentraland (identity swap with row 45)
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.6227"
$ws.Range("E46").Value = "  -1.17%  "

# Row 47 - NEARProtocol -> EnergySwap (identity swap with row 48)
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'13.42"
$ws.Range("E47").Value = "  -0.11%  "

# Row 48 - EnergySwap -> NEARProtocol (identity swap with row 47)
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'2.196"
$ws.Range("E48").Value = "  +0.15%  "

# Row 49 - PancakeSwap
$ws.Range("D49").Value = "'3.484"
$ws.Range("E49").Value = "  -2.02%  "

# Row 50 - BabyDogeCoin
$ws.Range("E50").Value = "  +0.14%  "

# Row 51 - Cronos
$ws.Range("D51").Value = "'0.06975"
$ws.Range("E51").Value = "  -0.19%  "
